# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de report sheets, as part of regenerating the
# Handback status report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 timestamps move forward
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 08:46:51"
$wsZhCn.Range("H2").Value = "2016-03-14 08:47:13"

# de-de sheet: row 2 timestamps move forward
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 08:46:55"
$wsDeDe.Range("H2").Value = "2016-03-14 08:47:22"
